$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the threshold value in C2 (11 -> 10.8)
$ws.Range("C2").Value = 10.8

# Move the active selection to C2 to match the saved view state
$ws.Range("C2").Select()
